$wb = $excel.ActiveWorkbook

# --- Create the new "Insertion Sort" sheet by copying the structurally
# identical "Selection Sort" sheet (same row layout for the best/worst
# case summary blocks), then renaming it. ---
$src = $wb.Worksheets.Item("Selection Sort")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "Insertion Sort"

# --- Header / title text ---
$new.Range("B2").Value = "Ratio Analysis of Insertion Sorting"
$new.Range("D3").Value = "Number of overwritings"

# --- Main data table (rows 4-11): columns C (count), D (overwritings),
# E (execution time) are literal values; F:O recompute automatically
# since they are formulas referencing B/C/F/G/H/I/J. ---
$new.Range("D4").Value = 0
$new.Range("E4").Value = 0.0000011600000000000399

$new.Range("C5").Value = 19.399999999999999
$new.Range("D5").Value = 19.399999999999999
$new.Range("E5").Value = 0.00000431999999999932

$new.Range("C6").Value = 663
$new.Range("D6").Value = 663
$new.Range("E6").Value = 0.000087160000000002696

$new.Range("C7").Value = 2543
$new.Range("D7").Value = 2543
$new.Range("E7").Value = 0.00026323999999999499

$new.Range("C8").Value = 62857.599999999999
$new.Range("D8").Value = 62857.599999999999
$new.Range("E8").Value = 0.0061542999999999997

$new.Range("C9").Value = 252100.2
$new.Range("D9").Value = 252100.2
$new.Range("E9").Value = 0.026573079999999999

$new.Range("C10").Value = 6222192.4000000004
$new.Range("D10").Value = 6222192.4000000004
$new.Range("E10").Value = 0.65331218000000002

$new.Range("C11").Value = 25013194.800000001
$new.Range("D11").Value = 25013194.800000001
$new.Range("E11").Value = 2.6780067599999899

# --- Best Case summary block (rows 14-17) ---
$new.Range("D15").Value = 0
$new.Range("D16").Value = 0
$new.Range("D17").Value = 0.00000299999999953115

# --- Worst Case summary block (rows 19-22) ---
$new.Range("D21").Value = 45
$new.Range("D22").Value = 0.0000056999999991091903

# --- Sheet view: selection on D5, this sheet becomes the active tab ---
$new.Range("D5").Select()
